# Apply "UserStories korrigiert von Marvin" corrections to the Flattie user
# stories sheet: update the wording of several user-story cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B10").Value = "As a user, I want to be able to throw myself out of a Flattie group, so that I can join or create another group. "
$ws.Range("B15").Value = "As a user, I want to be able to add a repeatable event into the shared calendar, so that I don't have to enter same events multiple times manually."
$ws.Range("B17").Value = "As a user, I want to be able to delete a repeatable event, so that I don't have to delete same events multiple times manually."
$ws.Range("B20").Value = "As a user, I want to be able to update a repeatable event, so that I don't have to change the same event information multiple times manually."
$ws.Range("B21").Value = "As a user, I want to be able to update one entry of a repeatable event, so that I can enter exceptions for the repeatable event."
$ws.Range("B23").Value = "As a user, I want to be able to assign a calendar event to a Flattie group member, so that all Flattie group members can see who is responsible for the event (especially important for duties)."
$ws.Range("B26").Value = "As a user, I want to be able to add an item to the shopping list, so that the needed item will get bought by a member of the group."

# Match the final view state recorded in the workbook (B30 selected).
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("B30").Select()

$wb.Save()

